$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "43.027.47"
$ws.Cells.Item(2, 5).Value = "  +2.11%  "

$ws.Cells.Item(3, 4).Value = "2.309.68"
$ws.Cells.Item(3, 5).Value = "  +1.82%  "

$ws.Cells.Item(4, 5).Value = "  +0.00%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "303.72"
$ws.Cells.Item(5, 5).Value = "  +2.10%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "100.76"
$ws.Cells.Item(6, 5).Value = "  +6.39%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.505"
$ws.Cells.Item(7, 5).Value = "  +2.94%  "

$ws.Cells.Item(8, 5).Value = "  -0.01%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.509"
$ws.Cells.Item(9, 5).Value = "  +3.66%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "34.86"
$ws.Cells.Item(10, 5).Value = "  +4.59%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0797"
$ws.Cells.Item(11, 5).Value = "  +1.05%  "

$ws.Cells.Item(12, 5).Value = "  +4.13%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "17.93"
$ws.Cells.Item(13, 5).Value = "  +14.06%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.86"
$ws.Cells.Item(14, 5).Value = "  +3.23%  "

$ws.Cells.Item(15, 4).Value = "2.687.82"
$ws.Cells.Item(15, 5).Value = "  +2.59%  "

$ws.Cells.Item(16, 4).Value = "2.299.86"
$ws.Cells.Item(16, 5).Value = "  +1.48%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.815"
$ws.Cells.Item(17, 5).Value = "  +4.53%  "

$ws.Cells.Item(18, 4).Value = "42.962.45"
$ws.Cells.Item(18, 5).Value = "  +2.10%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "12.48"
$ws.Cells.Item(19, 5).Value = "  +7.65%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "6.16"
$ws.Cells.Item(20, 5).Value = "  +3.01%  "

$ws.Cells.Item(21, 4).Value = "0.0₃0905"
$ws.Cells.Item(21, 5).Value = "  +1.83%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "67.77"
$ws.Cells.Item(22, 5).Value = "  +1.98%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "236.89"
$ws.Cells.Item(23, 5).Value = "  +1.59%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.22"
$ws.Cells.Item(24, 5).Value = "  +12.85%  "

$ws.Cells.Item(25, 2).Value = "Dai"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "1.00"
$ws.Cells.Item(25, 5).Value = "  -0.10%  "

$ws.Cells.Item(26, 2).Value = "PancakeSwap"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.47"
$ws.Cells.Item(26, 5).Value = "  +0.87%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "24.79"
$ws.Cells.Item(27, 5).Value = "  +3.88%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "168.28"
$ws.Cells.Item(28, 5).Value = "  +0.36%  "

$ws.Cells.Item(29, 5).Value = "  -7.90%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "34.07"
$ws.Cells.Item(30, 5).Value = "  +0.19%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "9.16"
$ws.Cells.Item(31, 5).Value = "  +0.86%  "

$ws.Cells.Item(32, 5).Value = "  +0.02%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "5.03"
$ws.Cells.Item(33, 5).Value = "  +2.86%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "4.63"
$ws.Cells.Item(34, 5).Value = "  +3.16%  "

$ws.Cells.Item(35, 2).Value = "WEMIXToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.42"
$ws.Cells.Item(35, 5).Value = "  +4.22%  "

$ws.Cells.Item(36, 2).Value = "Celestia"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "17.11"
$ws.Cells.Item(36, 5).Value = "  +4.05%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.0693"
$ws.Cells.Item(37, 5).Value = "  +0.79%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.102"
$ws.Cells.Item(38, 5).Value = "  +3.82%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.79"
$ws.Cells.Item(39, 5).Value = "  +4.20%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.82"
$ws.Cells.Item(40, 5).Value = "  +1.68%  "

$ws.Cells.Item(41, 5).Value = "  +0.78%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "2.33"
$ws.Cells.Item(42, 5).Value = "  -4.40%  "

$ws.Cells.Item(43, 4).Value = "2.002.85"
$ws.Cells.Item(43, 5).Value = "  +2.22%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0287"
$ws.Cells.Item(44, 5).Value = "  +3.30%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "10.16"
$ws.Cells.Item(45, 5).Value = "  +6.43%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "17.56"
$ws.Cells.Item(46, 5).Value = "  +1.15%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.85"
$ws.Cells.Item(47, 5).Value = "  +2.41%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "55.79"
$ws.Cells.Item(48, 5).Value = "  +7.55%  "

$ws.Cells.Item(49, 2).Value = "RocketPoolETH"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(49, 4).Value = "2.531.04"
$ws.Cells.Item(49, 5).Value = "  +1.51%  "

$ws.Cells.Item(50, 2).Value = "Stacks"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.54"
$ws.Cells.Item(50, 5).Value = "  +5.73%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "4.57"
$ws.Cells.Item(51, 5).Value = "  +0.77%  "
